# Insert a new data row at row 192 (pushing existing rows 192..262 down to 193..263)
# and populate it with the new observation's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows down by inserting a blank row at 192.
$ws.Rows.Item(192).Insert()

# Populate the newly inserted row 192 with the new record's data.
$ws.Cells.Item(192, 1).Value = 10
$ws.Cells.Item(192, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(192, 3).Value = "La Araucanía"
$ws.Cells.Item(192, 4).Value = 44489
$ws.Cells.Item(192, 5).Value = 9
$ws.Cells.Item(192, 6).Value = 100112032
$ws.Cells.Item(192, 7).Value = "Zapallo italiano"
$ws.Cells.Item(192, 8).Value = "Sin especificar"
$ws.Cells.Item(192, 9).Value = "Primera"
$ws.Cells.Item(192, 10).Value = 175
$ws.Cells.Item(192, 11).Value = 12000
$ws.Cells.Item(192, 12).Value = 13000
$ws.Cells.Item(192, 13).Value = 12457
$ws.Cells.Item(192, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(192, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(192, 16).Value = 208
$ws.Cells.Item(192, 17).Value = 60
$ws.Cells.Item(192, 18).Value = "Hortaliza"
